$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: add comment on row 21 (new F21 cell) ---
$ws1.Range("F21").Value = "Straight forward after doing so much of it"

# --- Sheet1: row 26 edits (hours bumped, comment extended) ---
$ws1.Range("B26").Value = 0.5
$ws1.Range("F26").Value = "Once you’ve done one…plus bonus dynamic selectors"

# --- Sheet1: new row 27 (Bugfix - amount_int to float) ---
$ws1.Range("A27").Value = 43505
$ws1.Range("A27").NumberFormat = $ws1.Range("A26").NumberFormat
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = "Bugfix"
$ws1.Range("D27").Value = 1
$ws1.Range("E27").Value = "Transaction.amount_int to float"
$ws1.Range("F27").Value = "Well. This is embarrassing. D*ck. Monetary values should be float. Ground up howler."
$ws1.Rows.Item(27).RowHeight = 16.65

# --- Sheet1: new row 28 (Bugfix - trailing zero truncation) ---
$ws1.Range("A28").Value = 43505
$ws1.Range("A28").NumberFormat = $ws1.Range("A26").NumberFormat
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = "Bugfix"
$ws1.Range("D28").Value = 1
$ws1.Range("E28").Value = "Transaction total truncating trailing zero"
$ws1.Range("F28").Value = "Fixed by converting to string using 2f string formatting"

# --- Sheet1: new row 29 (Planning - Cycle 2) ---
$ws1.Range("A29").Value = 43505
$ws1.Range("A29").NumberFormat = $ws1.Range("A26").NumberFormat
$ws1.Range("C29").Value = "Planning"
$ws1.Range("D29").Value = 2
$ws1.Range("E29").Value = "Discover user needs and measure against time remaining"
$ws1.Range("F29").Value = "Have to be careful, play it safe – leave Tue/Wed for CSS malark."

# --- Sheet2: view state (selection touches B7:B20, alongside the pre-existing B27) ---
$ws2.Select()
$ws2.Range("B7:B20").Select()

# --- Sheet1: view state (scrolled down, B7:B20 selected) - restore as active sheet ---
$ws1.Select()
$ws1.Range("B7:B20").Select()
